# Append 2023 and 2024 monthly "Total Fertilizer Production" figures to the
# existing time series on Sheet1 (columns A = month-start date, B = value).
# Mirrors the author's manual update described in the commit message:
# "updated total fertilizer production to 2023 - manual"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel date serials for the first of each month, Jan-2023 .. Dec-2023
$dates2023 = @(44927, 44958, 44986, 45017, 45047, 45078, 45108, 45139, 45170, 45200, 45231, 45261)
$val2023 = 289233334.27808142

# Excel date serials for the first of each month, Jan-2024 .. Dec-2024
$dates2024 = @(45292, 45323, 45352, 45383, 45413, 45444, 45474, 45505, 45536, 45566, 45597, 45627)
$val2024 = 297817126.50532985

# Existing data runs through row 181 (last date 2022-12-01); new rows continue at 182.
$row = 182

foreach ($d in $dates2023) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 1).NumberFormat = "mmm-yy"
    $ws.Cells.Item($row, 2).Value = $val2023
    $row = $row + 1
}

foreach ($d in $dates2024) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 1).NumberFormat = "mmm-yy"
    $ws.Cells.Item($row, 2).Value = $val2024
    $row = $row + 1
}

# Match the author's final view state: scrolled down, new rows selected.
$ws.Range("A172").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 172
$win.ScrollColumn = 1
$ws.Range("B194:B205").Select()
